# Weekly update: insert two new daily records (2023-12-07, serial 45267)
# at the top of the data block (row 1057), pushing all existing records
# down by two rows (old 1057-1116 -> new 1059-1118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 1057 (old first data-only row
# of this date group), shifting rows 1057:1116 down to 1059:1118.
$ws.Range("1057:1058").Insert()

# ---- New row 1057 : Coliflor, Primera ----
$ws.Cells.Item(1057, 1).Value = 9
$ws.Cells.Item(1057, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1057, 3).Value = 'Metropolitana'
$ws.Cells.Item(1057, 4).Value = 45267
$ws.Cells.Item(1057, 5).Value = 13
$ws.Cells.Item(1057, 6).Value = 100112008
$ws.Cells.Item(1057, 7).Value = 'Coliflor'
$ws.Cells.Item(1057, 8).Value = 'Sin especificar'
$ws.Cells.Item(1057, 9).Value = 'Primera'
$ws.Cells.Item(1057, 10).Value = 1600
$ws.Cells.Item(1057, 11).Value = 800
$ws.Cells.Item(1057, 12).Value = 900
$ws.Cells.Item(1057, 13).Value = 850
$ws.Cells.Item(1057, 14).Value = '$/unidad'
$ws.Cells.Item(1057, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(1057, 16).Value = 850
$ws.Cells.Item(1057, 17).Value = 1
$ws.Cells.Item(1057, 18).Value = 'Hortaliza'

# ---- New row 1058 : Coliflor, Segunda ----
$ws.Cells.Item(1058, 1).Value = 9
$ws.Cells.Item(1058, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1058, 3).Value = 'Metropolitana'
$ws.Cells.Item(1058, 4).Value = 45267
$ws.Cells.Item(1058, 5).Value = 13
$ws.Cells.Item(1058, 6).Value = 100112008
$ws.Cells.Item(1058, 7).Value = 'Coliflor'
$ws.Cells.Item(1058, 8).Value = 'Sin especificar'
$ws.Cells.Item(1058, 9).Value = 'Segunda'
$ws.Cells.Item(1058, 10).Value = 970
$ws.Cells.Item(1058, 11).Value = 700
$ws.Cells.Item(1058, 12).Value = 700
$ws.Cells.Item(1058, 13).Value = 700
$ws.Cells.Item(1058, 14).Value = '$/unidad'
$ws.Cells.Item(1058, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(1058, 16).Value = 700
$ws.Cells.Item(1058, 17).Value = 1
$ws.Cells.Item(1058, 18).Value = 'Hortaliza'
